# Update NATMI TPM-derived statistics for the Fgf17-Fgfr1 LR-pair sheet
# with newly recomputed TPM values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.6594814792829158
$ws.Range("J2").Value = 0.6594814792829158
$ws.Range("M2").Value = 1.845768666666667
$ws.Range("N2").Value = 5.537306
$ws.Range("O2").Value = 0.01459089321241885
$ws.Range("P2").Value = 0.01459089321241885
$ws.Range("Q2").Value = 0.1093642545248889
$ws.Range("R2").Value = 0.984278290724
$ws.Range("S2").Value = 0.009622423839785038
$ws.Range("T2").Value = 0.009622423839785038

# Row 3
$ws.Range("I3").Value = 0.6594814792829158
$ws.Range("J3").Value = 0.6594814792829158
$ws.Range("O3").Value = 0.6557810310272387
$ws.Range("P3").Value = 0.6557810310272387
$ws.Range("S3").Value = 0.4324754444275191
$ws.Range("T3").Value = 0.4324754444275191

# Row 4
$ws.Range("I4").Value = 0.6594814792829158
$ws.Range("J4").Value = 0.6594814792829158
$ws.Range("M4").Value = 41.69841866666667
$ws.Range("N4").Value = 125.095256
$ws.Range("O4").Value = 0.3296280757603424
$ws.Range("P4").Value = 0.3296280757603424
$ws.Range("Q4").Value = 2.470686903891556
$ws.Range("R4").Value = 22.236182135024
$ws.Range("S4").Value = 0.2173836110156117
$ws.Range("T4").Value = 0.2173836110156117

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.030594
$ws.Range("H5").Value = 0.091782
$ws.Range("I5").Value = 0.3405185207170842
$ws.Range("J5").Value = 0.3405185207170842
$ws.Range("M5").Value = 1.845768666666667
$ws.Range("N5").Value = 5.537306
$ws.Range("O5").Value = 0.01459089321241885
$ws.Range("P5").Value = 0.01459089321241885
$ws.Range("Q5").Value = 0.056469446588
$ws.Range("R5").Value = 0.508225019292
$ws.Range("S5").Value = 0.004968469372633811
$ws.Range("T5").Value = 0.004968469372633811

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.030594
$ws.Range("H6").Value = 0.091782
$ws.Range("I6").Value = 0.3405185207170842
$ws.Range("J6").Value = 0.3405185207170842
$ws.Range("O6").Value = 0.6557810310272387
$ws.Range("P6").Value = 0.6557810310272387
$ws.Range("Q6").Value = 2.537993484422
$ws.Range("R6").Value = 22.841941359798
$ws.Range("S6").Value = 0.2233055865997196
$ws.Range("T6").Value = 0.2233055865997196

# Row 7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.030594
$ws.Range("H7").Value = 0.091782
$ws.Range("I7").Value = 0.3405185207170842
$ws.Range("J7").Value = 0.3405185207170842
$ws.Range("M7").Value = 41.69841866666667
$ws.Range("N7").Value = 125.095256
$ws.Range("O7").Value = 0.3296280757603424
$ws.Range("P7").Value = 0.3296280757603424
$ws.Range("Q7").Value = 1.275721420688
$ws.Range("R7").Value = 11.481492786192
$ws.Range("S7").Value = 0.1122444647447308
$ws.Range("T7").Value = 0.1122444647447308

Write-Host "Updated TPM values for rows 2-7"
